# Reflect the new requirement separation on the ITP sheet:
#   - Old column D ("Terms Typically Offered") moves to column G.
#   - Three new columns are inserted at D/E/F: Corequisites, Concurrent,
#     Recommended. Every existing row gets "NA" in these new columns,
#     except row 8 (ITP 275) whose corequisite text ("ITP 211.") moves out
#     of the Prerequisites cell (C8) into the new Corequisites cell (D8).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Terms-typically-offered values, keyed by row, that will land in column G.
# These are simply the current column-D values, carried over verbatim,
# except row 8 which becomes "SP " (trailing space) per the source data.
$termsByRow = @{
    2 = 'F, W, SP'
    3 = 'F, W'
    4 = 'F, W, SP'
    5 = 'F'
    6 = 'W, SP'
    7 = 'TBD'
    8 = 'SP '
    9 = 'W'
    10 = 'F,W,SP,SU'
    11 = 'F, SP'
    12 = 'F,W,SP,SU'
    13 = 'F, W, SP'
    14 = 'F,W,SP,SU'
    15 = 'F, W'
    16 = 'TBD'
    17 = 'F, W'
    18 = 'SP'
    19 = 'F, SP'
    20 = 'W'
    21 = 'W'
    22 = 'F, W'
    23 = 'F'
    24 = 'W'
    25 = 'SP'
    26 = 'W'
    27 = 'SP'
    28 = 'F, W, SP'
    29 = 'F, W, SP'
    30 = 'W'
    31 = 'TBD'
    32 = 'F'
    33 = 'W'
    34 = 'TBD'
    35 = 'F, SP'
    36 = 'TBD'
    37 = 'TBD'
    38 = 'SP'
    39 = 'SP'
    40 = 'TBD'
    41 = 'TBD'
    42 = 'TBD'
    43 = 'TBD'
    44 = 'TBD'
    45 = 'TBD'
    46 = 'TBD'
}

# Header row: insert the three new headers before "Terms Typically Offered"
# which shifts from D1 to G1.
$ws.Cells.Item(1, 4).Value = "Corequisites"
$ws.Cells.Item(1, 5).Value = "Concurrent"
$ws.Cells.Item(1, 6).Value = "Recommended"
$ws.Cells.Item(1, 7).Value = "Terms Typically Offered"

# Row 8 (ITP 275) calls out an actual corequisite, so split it out of the
# Prerequisites text and into the new Corequisites column.
$ws.Cells.Item(8, 3).Value = "ITP 150."
$ws.Cells.Item(8, 4).Value = "ITP 211."

# Two prerequisite-text cleanups that accompany the column restructuring.
$ws.Cells.Item(14, 3).Value = "A grade of C- or better, or consent of instructor, MATH 141 or MATH 221, and STAT 217 or STAT 218 or STAT 251 or any 300 or 400 level statistics course."
$ws.Cells.Item(27, 3).Value = "One of the ITP 303, ITP 326, ITP 330, or ITP 341; and ITP 371."

# Data rows 2-46: fill Corequisites / Concurrent / Recommended with "NA"
# (row 8's Corequisites cell was already set above) and move the old
# Terms-Typically-Offered value into column G.
for ($row = 2; $row -le 46; $row++) {
    if ($row -ne 8) {
        $ws.Cells.Item($row, 4).Value = "NA"
    }
    $ws.Cells.Item($row, 5).Value = "NA"
    $ws.Cells.Item($row, 6).Value = "NA"
    $ws.Cells.Item($row, 7).Value = $termsByRow[$row]
}
